# Regenerate the localization status report: the four files that were
# "In Translation" / "Ready for handoff" are re-ordered (2a751b94 and
# 625926c8 move to "In Translation", 88989a8a stays "In Translation" but
# moves down, 334ba757 remains "Ready for handoff" last) across the
# Overview, zh-cn and de-de sheets (rows 7-10).

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $addr = $cell.Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $value
        }
    }
}

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @{
    7  = @{ A = "2a751b94-b719-4f3b-b96f-c6fe889bff3c.md"; B = "In Translation";    C = "In Translation";    D = "2016-17-11 22:17:06" }
    8  = @{ A = "625926c8-f63d-4c5d-9c3b-90a800cdb0ce.md"; B = "In Translation";    C = "In Translation";    D = "2016-17-11 22:17:06" }
    9  = @{ A = "88989a8a-0262-4447-ab47-0b581d8ecc69.md"; B = "In Translation";    C = "In Translation";    D = "2016-12-11 22:12:26" }
    10 = @{ A = "334ba757-3e04-4960-b90c-ec9a3eb937d1.md"; B = "Ready for handoff"; C = "Ready for handoff"; D = "2016-17-11 22:17:06" }
}

foreach ($r in 7..10) {
    $row = $overviewRows[$r]
    Set-CellAndHyperlink $wsOverview $r 1 $row.A
    $wsOverview.Cells.Item($r, 2).Value = $row.B
    $wsOverview.Cells.Item($r, 3).Value = $row.C
    $wsOverview.Cells.Item($r, 4).Value = $row.D
}

# ---------------- zh-cn sheet ----------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnRows = @{
    7  = @{ A = "2a751b94-b719-4f3b-b96f-c6fe889bff3c.md"; B = ".md"; C = "In Translation";    D = "2a751b94-b719-4f3b-b96f-c6fe889bff3c.5294792b89841056768df96fa2cb2bbd8cfa559c.zh-cn.xlf"; E = "2016-03-11 22:17:03"; H = "0001-01-01 00:00:00"; I = "Include" }
    8  = @{ A = "625926c8-f63d-4c5d-9c3b-90a800cdb0ce.md"; B = ".md"; C = "In Translation";    D = "625926c8-f63d-4c5d-9c3b-90a800cdb0ce.17ce36c2ac489c2cd7280253c6803abaf14f9bef.zh-cn.xlf"; E = "2016-03-11 22:17:03"; H = "0001-01-01 00:00:00"; I = "Include" }
    9  = @{ A = "88989a8a-0262-4447-ab47-0b581d8ecc69.md"; B = ".md"; C = "In Translation";    D = "88989a8a-0262-4447-ab47-0b581d8ecc69.8161ab9136248d51b8caa3fafd358546e9cb76a5.zh-cn.xlf"; E = "2016-03-11 22:12:16"; H = "0001-01-01 00:00:00"; I = "Include" }
    10 = @{ A = "334ba757-3e04-4960-b90c-ec9a3eb937d1.md"; B = ".md"; C = "Ready for handoff"; D = "334ba757-3e04-4960-b90c-ec9a3eb937d1.f1c51fedfdf63615d256bce43d1c1a62ece4c6bb.zh-cn.xlf"; E = "2016-03-11 22:17:03"; H = "0001-01-01 00:00:00"; I = "Include" }
}

foreach ($r in 7..10) {
    $row = $zhCnRows[$r]
    Set-CellAndHyperlink $wsZhCn $r 1 $row.A
    Set-CellAndHyperlink $wsZhCn $r 2 $row.B
    $wsZhCn.Cells.Item($r, 3).Value = $row.C
    Set-CellAndHyperlink $wsZhCn $r 4 $row.D
    $wsZhCn.Cells.Item($r, 5).Value = $row.E
    $wsZhCn.Cells.Item($r, 8).Value = $row.H
    $wsZhCn.Cells.Item($r, 9).Value = $row.I
}

# ---------------- de-de sheet ----------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeRows = @{
    7  = @{ A = "2a751b94-b719-4f3b-b96f-c6fe889bff3c.md"; B = ".md"; C = "In Translation";    D = "2a751b94-b719-4f3b-b96f-c6fe889bff3c.5294792b89841056768df96fa2cb2bbd8cfa559c.de-de.xlf"; E = "2016-03-11 22:17:06"; H = "0001-01-01 00:00:00"; I = "Include" }
    8  = @{ A = "625926c8-f63d-4c5d-9c3b-90a800cdb0ce.md"; B = ".md"; C = "In Translation";    D = "625926c8-f63d-4c5d-9c3b-90a800cdb0ce.17ce36c2ac489c2cd7280253c6803abaf14f9bef.de-de.xlf"; E = "2016-03-11 22:17:06"; H = "0001-01-01 00:00:00"; I = "Include" }
    9  = @{ A = "88989a8a-0262-4447-ab47-0b581d8ecc69.md"; B = ".md"; C = "In Translation";    D = "88989a8a-0262-4447-ab47-0b581d8ecc69.8161ab9136248d51b8caa3fafd358546e9cb76a5.de-de.xlf"; E = "2016-03-11 22:12:26"; H = "0001-01-01 00:00:00"; I = "Include" }
    10 = @{ A = "334ba757-3e04-4960-b90c-ec9a3eb937d1.md"; B = ".md"; C = "Ready for handoff"; D = "334ba757-3e04-4960-b90c-ec9a3eb937d1.f1c51fedfdf63615d256bce43d1c1a62ece4c6bb.de-de.xlf"; E = "2016-03-11 22:17:06"; H = "0001-01-01 00:00:00"; I = "Include" }
}

foreach ($r in 7..10) {
    $row = $deDeRows[$r]
    Set-CellAndHyperlink $wsDeDe $r 1 $row.A
    Set-CellAndHyperlink $wsDeDe $r 2 $row.B
    $wsDeDe.Cells.Item($r, 3).Value = $row.C
    Set-CellAndHyperlink $wsDeDe $r 4 $row.D
    $wsDeDe.Cells.Item($r, 5).Value = $row.E
    $wsDeDe.Cells.Item($r, 8).Value = $row.H
    $wsDeDe.Cells.Item($r, 9).Value = $row.I
}
